# Updates "北京-漫展信息.xlsx" to the figures captured at commit 456a3b4.
# Touches three of the four sheets:
#   - 展览   (sheet1) : refreshed "want to go" / price counters, new cover image for row 48
#   - 演出   (sheet2) : refreshed counters + one brand-new event row (row 10)
#   - 全部类型 (sheet4): the merged view, refreshed with the same counters
# 本地生活 (sheet3) is untouched in this revision.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F4").Value  = 447
$ws1.Range("F5").Value  = 1314
$ws1.Range("F6").Value  = 7643
$ws1.Range("F9").Value  = 2090
$ws1.Range("F10").Value = 8440
$ws1.Range("F14").Value = 5652
$ws1.Range("F16").Value = 2608
$ws1.Range("F17").Value = 1131
$ws1.Range("F22").Value = 32
$ws1.Range("G22").Value = 35.1
$ws1.Range("F23").Value = 528
$ws1.Range("F24").Value = 3487
$ws1.Range("F25").Value = 44
$ws1.Range("F27").Value = 20
$ws1.Range("F29").Value = 2976
$ws1.Range("F30").Value = 17
$ws1.Range("F31").Value = 82
$ws1.Range("F32").Value = 340
$ws1.Range("F34").Value = 307
$ws1.Range("F35").Value = 304
$ws1.Range("F36").Value = 654
$ws1.Range("F39").Value = 1721
$ws1.Range("F42").Value = 16
$ws1.Range("F43").Value = 2749
$ws1.Range("F45").Value = 2283
$ws1.Range("F46").Value = 9
$ws1.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202403/dAYut7iv1709883417318.jpeg"

# ---------------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F3").Value = 121
$ws2.Range("F4").Value = 6
$ws2.Range("F9").Value = 115

# New row 10 - new event that appeared in this revision.
# Copy A9's formatting (bold/centered/bordered index-column style) onto A10,
# then overwrite with the real value, so the new row matches the look of
# every other row's leftmost "index" column.
$ws2.Range("A9").Copy($ws2.Range("A10"))
$ws2.Range("A10").Value = 9
# Force text so the date-looking string isn't silently coerced into a date serial.
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "2024-06-28"
$ws2.Range("C10").Value = "北京·《国风大赏》大型国潮音乐会×郑州歌舞剧院《唐宫夜宴》"
$ws2.Range("D10").Value = "西直门外大街135号（北京展览馆内） 北京展览馆剧场"
$ws2.Range("E10").Value = "2024.06.28 19:30-06.28 21:00"
$ws2.Range("F10").Value = 0
$ws2.Range("G10").Value = 162
$ws2.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=82587"
$ws2.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202403/VZcJ2SJ51709882503997.jpeg"

# ---------------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value  = 1314
$ws4.Range("F6").Value  = 7643
$ws4.Range("F9").Value  = 2090
$ws4.Range("F10").Value = 8440
$ws4.Range("F14").Value = 5652
$ws4.Range("F16").Value = 2608
$ws4.Range("F17").Value = 1131
$ws4.Range("F22").Value = 32
$ws4.Range("G22").Value = 35.1
$ws4.Range("F23").Value = 121
$ws4.Range("F24").Value = 528
$ws4.Range("F25").Value = 6
$ws4.Range("F26").Value = 3487
$ws4.Range("F28").Value = 20
$ws4.Range("F29").Value = 2976
$ws4.Range("F30").Value = 340
$ws4.Range("F32").Value = 307
$ws4.Range("F34").Value = 304
$ws4.Range("F35").Value = 654
$ws4.Range("F39").Value = 1721
$ws4.Range("F42").Value = 16
$ws4.Range("F43").Value = 2751
$ws4.Range("F46").Value = 2283
$ws4.Range("I48").Value = "//i2.hdslb.com/bfs/openplatform/202403/dAYut7iv1709883417318.jpeg"
$ws4.Range("F49").Value = 115
